# POV / Persona Korrekturen: kleinere Text- und Zeitangabe-Fixes.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                      $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found -> $old"
    }
}

# "...um 6:30 in sein Auto" -> "...um 6:30 Uhr in sein Auto"
Replace-Text "um 6:30 in sein Auto" "um 6:30 Uhr in sein Auto"

# "...weil man etwas kaufen muss, bestellt..." -> "...weil man dort etwas kaufen muss um einen Platz zu besetzen, bestellt..."
Replace-Text "weil man etwas kaufen muss," "weil man dort etwas kaufen muss um einen Platz zu besetzen,"

# "...dabei wird sich verabredet sich das Taxi..." -> "...dabei wird verabredet sich das Taxi..."
Replace-Text "dabei wird sich verabredet sich das Taxi" "dabei wird verabredet sich das Taxi"

# "...nach 16:00 ist geht keiner..." -> "...nach 16:00 Uhr ist geht keiner..."
Replace-Text "nach 16:00 ist geht keiner" "nach 16:00 Uhr ist geht keiner"

# "...die Mobilnummer des hauptverantwortlichen..." -> "...die Mobilfunknummer des hauptverantwortlichen..."
Replace-Text "Mobilnummer des hauptverantwortlichen" "Mobilfunknummer des hauptverantwortlichen"

# "...nachgucken. 10 Minuten später..." -> "...nachgucken. Zehn Minuten später..."
Replace-Text "nachgucken. 10 Minuten später" "nachgucken. Zehn Minuten später"

# "...auf Anhieb 3 kleiner Bugs..." -> "...auf Anhieb drei kleiner Bugs..."
Replace-Text "auf Anhieb 3 kleiner Bugs" "auf Anhieb drei kleiner Bugs"

# "...um 18.18 erreichen..." -> "...um 18:18 Uhr erreichen..."
Replace-Text "um 18.18 erreichen" "um 18:18 Uhr erreichen"
